$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.424719
$ws.Range("N2").Value = 4.274157
$ws.Range("O2").Value = 0.07423298812267187
$ws.Range("P2").Value = 0.07423298812267187
$ws.Range("Q2").Value = 0.06767557721899999
$ws.Range("R2").Value = 0.609080194971
$ws.Range("S2").Value = 0.07423298812267187
$ws.Range("T2").Value = 0.07423298812267187

# Row 3
$ws.Range("O3").Value = 0.5596266124066729
$ws.Range("P3").Value = 0.5596266124066729
$ws.Range("S3").Value = 0.5596266124066729
$ws.Range("T3").Value = 0.5596266124066729

# Row 4
$ws.Range("M4").Value = 7.027161333333335
$ws.Range("O4").Value = 0.3661403994706553
$ws.Range("P4").Value = 0.3661403994706552
$ws.Range("Q4").Value = 0.3337971904946667
$ws.Range("S4").Value = 0.3661403994706553
$ws.Range("T4").Value = 0.3661403994706552
